$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing data rows 2-14 (refreshed timestamps/hashes/fees/profit) ---
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = $true
$ws.Cells.Item(2, 3).Value = "Mon Apr 26 2021 01:22:53 GMT+0900 (Japan Standard Time)"
$ws.Cells.Item(2, 4).Value = "0x03d76e2dcf30c4095f529d905add2df604b7823babc59270f28f7d5b316edfc9"
$ws.Cells.Item(2, 5).Value = "0.000414524"
$ws.Cells.Item(2, 6).Value = 0.01160065113598243

$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = $true
$ws.Cells.Item(3, 3).Value = "Mon Apr 26 2021 19:36:29 GMT+0900 (Japan Standard Time)"
$ws.Cells.Item(3, 4).Value = "0x5e39157180c7266ddc648bc1ca7e735202e6ad42caa3f1aa381f4fcba7aef269"
$ws.Cells.Item(3, 5).Value = "0.00041387"
$ws.Cells.Item(3, 6).Value = 0.000425553863434797

$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = $true
$ws.Cells.Item(4, 3).Value = "Mon Apr 26 2021 20:18:25 GMT+0900 (Japan Standard Time)"
$ws.Cells.Item(4, 4).Value = "0xf556f6ab7055d48f9ebc5cc3f5c3ef9fbd7f102c83510ec197a9b05ba96f3d07"
$ws.Cells.Item(4, 5).Value = "0.000414568"
$ws.Cells.Item(4, 6).Value = 0.01053781198224355

$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(5, 2).Value = $true
$ws.Cells.Item(5, 3).Value = "Mon Apr 26 2021 20:40:21 GMT+0900 (Japan Standard Time)"
$ws.Cells.Item(5, 4).Value = "0x003e2e03654b5c888e9c9142fda0e61e178b5ee0acf3741c3a4f2f28cc0af6c1"
$ws.Cells.Item(5, 5).Value = "0.000414568"
$ws.Cells.Item(5, 6).Value = 0.009564243733416767

$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 2).Value = $true
$ws.Cells.Item(6, 3).Value = "Mon Apr 26 2021 20:42:09 GMT+0900 (Japan Standard Time)"
$ws.Cells.Item(6, 4).Value = "0xdecc9f5ea2d2027f5448db538b5ce52eeadf7531b48b7c3fde80ef8484420240"
$ws.Cells.Item(6, 5).Value = "0.000413894"
$ws.Cells.Item(6, 6).Value = 0.02138397072097727

$ws.Cells.Item(7, 1).Value = 5
$ws.Cells.Item(7, 2).Value = $true
$ws.Cells.Item(7, 3).Value = "Mon Apr 26 2021 21:06:33 GMT+0900 (Japan Standard Time)"
$ws.Cells.Item(7, 4).Value = "0x9dc151486cd3e88c6160f73c8f857ce59038e1e9bc8e4e25fcf8163eeff7e302"
$ws.Cells.Item(7, 5).Value = "0.00041387"
$ws.Cells.Item(7, 6).Value = 0.01939733196672238

$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = $true
$ws.Cells.Item(8, 3).Value = "Mon Apr 26 2021 21:08:32 GMT+0900 (Japan Standard Time)"
$ws.Cells.Item(8, 4).Value = "0x9dd2f6eed61bd4ffa60e782fefd9ee4ad8d8d80b68f29b15109c63647923e90c"
$ws.Cells.Item(8, 5).Value = "0.000414568"
$ws.Cells.Item(8, 6).Value = 0.005334091779495808

$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = $true
$ws.Cells.Item(9, 3).Value = "Mon Apr 26 2021 21:23:20 GMT+0900 (Japan Standard Time)"
$ws.Cells.Item(9, 4).Value = "0x583282fd9a2d74154237c30ac1a79fec24859322845f2a663c7f6f8b5d021cee"
$ws.Cells.Item(9, 5).Value = "0.000414568"
$ws.Cells.Item(9, 6).Value = 0.008856762837462144

$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = $true
$ws.Cells.Item(10, 3).Value = "Mon Apr 26 2021 21:26:21 GMT+0900 (Japan Standard Time)"
$ws.Cells.Item(10, 4).Value = "0x6c77c8ba69d2e7b99323accffa81ce2cf9ab13096418bea81520d22f2aefc09b"
$ws.Cells.Item(10, 5).Value = "0.000413894"
$ws.Cells.Item(10, 6).Value = 0.02233797328233905

$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = $true
$ws.Cells.Item(11, 3).Value = "Mon Apr 26 2021 21:32:30 GMT+0900 (Japan Standard Time)"
$ws.Cells.Item(11, 4).Value = "0xeb79ca09a138cc93c297e9a23787ffa6d418d22716a9c58267fdfe694f3bed8c"
$ws.Cells.Item(11, 5).Value = "0.000414592"
$ws.Cells.Item(11, 6).Value = -0.000246525737690668

$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 2).Value = $true
$ws.Cells.Item(12, 3).Value = "Mon Apr 26 2021 21:32:49 GMT+0900 (Japan Standard Time)"
$ws.Cells.Item(12, 4).Value = "0xb91d2c95711ecba9bfd313e26b91554e0610810a7047ab27903db2b39d5d146f"
$ws.Cells.Item(12, 5).Value = "0.000414592"
$ws.Cells.Item(12, 6).Value = -0.000399706484896156

$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 2).Value = $true
$ws.Cells.Item(13, 3).Value = "Mon Apr 26 2021 22:51:38 GMT+0900 (Japan Standard Time)"
$ws.Cells.Item(13, 4).Value = "0xe3cc495009957ff54fa4b5b64d20573261f37879674605a7755b44d47c3a6c71"
$ws.Cells.Item(13, 5).Value = "0.00621924"
$ws.Cells.Item(13, 6).Value = -0.004770505782486944

$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(14, 2).Value = $true
$ws.Cells.Item(14, 3).Value = "Tue Apr 27 2021 00:03:52 GMT+0900 (Japan Standard Time)"
$ws.Cells.Item(14, 4).Value = "0x8d6fe5bd2e979f36208af5c4248484439efb4926b3dc78451e05db3fed362b04"
$ws.Cells.Item(14, 5).Value = "0.001865772"
$ws.Cells.Item(14, 6).Value = 0.03646707030010118

# --- Append new rows 15-20 for new trade records ---
$ws.Cells.Item(15, 1).Value = 13
$ws.Cells.Item(15, 2).Value = $true
$ws.Cells.Item(15, 3).Value = "Tue Apr 27 2021 00:07:50 GMT+0900 (Japan Standard Time)"
$ws.Cells.Item(15, 4).Value = "0x917829dc338fb8ae0cd3470699de86cfc32e136910ddfec3c62008137d55f1e8"
$ws.Cells.Item(15, 5).Value = "0.001862415"
$ws.Cells.Item(15, 6).Value = 0.003592710552158496

$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = $true
$ws.Cells.Item(16, 3).Value = "Tue Apr 27 2021 00:08:17 GMT+0900 (Japan Standard Time)"
$ws.Cells.Item(16, 4).Value = "0x182529a4af3c89b41bc118885447dc687bd8c4a6a2f387c1c366caf5c88a00d7"
$ws.Cells.Item(16, 5).Value = "0.001862631"
$ws.Cells.Item(16, 6).Value = 0.001561784184273331

$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = $true
$ws.Cells.Item(17, 3).Value = "Tue Apr 27 2021 00:34:57 GMT+0900 (Japan Standard Time)"
$ws.Cells.Item(17, 4).Value = "0x738b3facd8c4b6cf6f5887f0666ee9b5a6ddead4d2866fc6436df050ea626311"
$ws.Cells.Item(17, 5).Value = "0.000620877"
$ws.Cells.Item(17, 6).Value = 0.02332403784313836

$ws.Cells.Item(18, 1).Value = 16
$ws.Cells.Item(18, 2).Value = $false
$ws.Cells.Item(18, 3).Value = "Tue Apr 27 2021 00:35:53 GMT+0900 (Japan Standard Time)"
$ws.Cells.Item(18, 4).Value = "0x4dc19f33d2e76a2407fc75d90a9133bcbee0c5f01859693c7d007e3231403b4f"
$ws.Cells.Item(18, 5).Value = "0.000787794"
$ws.Cells.Item(18, 6).Value = -0.000787794

$ws.Cells.Item(19, 1).Value = 17
$ws.Cells.Item(19, 2).Value = $true
$ws.Cells.Item(19, 3).Value = "Tue Apr 27 2021 00:36:29 GMT+0900 (Japan Standard Time)"
$ws.Cells.Item(19, 4).Value = "0x2ccea5cf3fd393213296044bd37118b41735650c3c9f9a28fd6349edce05e3c3"
$ws.Cells.Item(19, 5).Value = "0.000621924"
$ws.Cells.Item(19, 6).Value = 0.02639851150879546

$ws.Cells.Item(20, 1).Value = 18
$ws.Cells.Item(20, 2).Value = $true
$ws.Cells.Item(20, 3).Value = "Tue Apr 27 2021 00:37:53 GMT+0900 (Japan Standard Time)"
$ws.Cells.Item(20, 4).Value = "0x4891251aadf085682b5b207f7867d6896ca0dfbbefeb398af8d19a45d2674081"
$ws.Cells.Item(20, 5).Value = "0.000621888"
$ws.Cells.Item(20, 6).Value = 0.0205843811667368

# --- Carry the column-A number style (s="1") onto the new rows, like the existing rows ---
$ws.Range("A14").Copy() | Out-Null
$ws.Range("A15").PasteSpecial(-4122) | Out-Null
$ws.Range("A16").PasteSpecial(-4122) | Out-Null
$ws.Range("A17").PasteSpecial(-4122) | Out-Null
$ws.Range("A18").PasteSpecial(-4122) | Out-Null
$ws.Range("A19").PasteSpecial(-4122) | Out-Null
$ws.Range("A20").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1, 1).Select() | Out-Null
